$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$ok = $find.Execute("2022年7月1日  星期五", $true, $false, $false, $false, $false, $true, 1, $false, "2022年7月1日  星期五^p天晴，今天去北苑体育馆拿了个快递，中午的太阳太晒了，嘎达^p", 2)

$p13 = $d.Paragraphs.Item(13)
$p13.Range.Text = "今天天气好"

$p13 = $d.Paragraphs.Item(13)
$markPos = $p13.Range.End - 1
$d.Range($markPos, $markPos + 1).Delete()

$p13m = $d.Paragraphs.Item(13)
$p13m.Format.Alignment = 0
Write-Output "set Alignment done"
